$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.202.88"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.686.71"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +10.13%  "
$ws.Range("E9").Value = "  +4.62%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.925.50"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "1.689.52"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("E15").Value = "  +4.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.22"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "27.211.55"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.98"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("E23").Value = "  +4.91%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").Value = "1.549.54"
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.15"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.75"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.26"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "1.833.50"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.790"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.61"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0110"
$ws.Range("E49").Value = "  +3.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.11%  "
$ws.Range("E51").Value = "  +1.17%  "
